# ------------------------------------------------------------------
# Refined metadata to be additional tab
#
# 1) Update the "time_taken" (col F) timestamps on the existing
#    "data" sheet to the new re-query time.
# 2) Add a new "metadata" worksheet (after "data") describing the
#    PanelApp query that produced the data sheet.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the F column ("time_taken") timestamps on "data" ---
$timeTaken = @(
    "2021-10-05 14:33:10.206748",
    "2021-10-05 14:33:10.206756",
    "2021-10-05 14:33:10.206759",
    "2021-10-05 14:33:10.206761",
    "2021-10-05 14:33:10.206764",
    "2021-10-05 14:33:10.206767",
    "2021-10-05 14:33:10.206769",
    "2021-10-05 14:33:10.206772",
    "2021-10-05 14:33:10.206775",
    "2021-10-05 14:33:10.206777",
    "2021-10-05 14:33:10.206780",
    "2021-10-05 14:33:10.206782",
    "2021-10-05 14:33:10.206785",
    "2021-10-05 14:33:10.206787",
    "2021-10-05 14:33:10.206790",
    "2021-10-05 14:33:10.206792",
    "2021-10-05 14:33:10.206795",
    "2021-10-05 14:33:10.206798",
    "2021-10-05 14:33:10.206801",
    "2021-10-05 14:33:10.206803",
    "2021-10-05 14:33:10.206806",
    "2021-10-05 14:33:10.206808",
    "2021-10-05 14:33:10.206811",
    "2021-10-05 14:33:10.206813",
    "2021-10-05 14:33:10.206816",
    "2021-10-05 14:33:10.206819",
    "2021-10-05 14:33:10.206822",
    "2021-10-05 14:33:10.206824",
    "2021-10-05 14:33:10.206827",
    "2021-10-05 14:33:10.206829",
    "2021-10-05 14:33:10.206832",
    "2021-10-05 14:33:10.206834",
    "2021-10-05 14:33:10.206837",
    "2021-10-05 14:33:10.206839",
    "2021-10-05 14:33:10.206842",
    "2021-10-05 14:33:10.206844",
    "2021-10-05 14:33:10.206847",
    "2021-10-05 14:33:10.206849",
    "2021-10-05 14:33:10.206852",
    "2021-10-05 14:33:10.206854",
    "2021-10-05 14:33:10.206857",
    "2021-10-05 14:33:10.206860",
    "2021-10-05 14:33:10.206862",
    "2021-10-05 14:33:10.206865",
    "2021-10-05 14:33:10.206867",
    "2021-10-05 14:33:10.206870",
    "2021-10-05 14:33:10.206872",
    "2021-10-05 14:33:10.206874",
    "2021-10-05 14:33:10.206877",
    "2021-10-05 14:33:10.206879",
    "2021-10-05 14:33:10.206882",
    "2021-10-05 14:33:10.206884",
    "2021-10-05 14:33:10.206887",
    "2021-10-05 14:33:10.206890",
    "2021-10-05 14:33:10.206892",
    "2021-10-05 14:33:10.206895",
    "2021-10-05 14:33:10.206897",
    "2021-10-05 14:33:10.206900",
    "2021-10-05 14:33:10.206902",
    "2021-10-05 14:33:10.206905",
    "2021-10-05 14:33:10.206907",
    "2021-10-05 14:33:10.206910",
    "2021-10-05 14:33:10.206913"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timeTaken[$i]
}

# --- 2. Add the new "metadata" worksheet, placed after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (row 1), columns B:G
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # B = 2
    $metaSheet.Cells.Item(1, $col).Value = $headers[$i]
}

# Copy the header cell formatting from "data"!B1:F1 (bold / bordered /
# centered style already present in the workbook's style table) onto
# the new header cells, re-using the existing style instead of
# creating a new one.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (row 2)
$metaSheet.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Range("B2").Value = "Arrhythmia_SuperPanel"
$metaSheet.Range("C2").Value = 254

# data_version ("1.2") must be stored as text, not as the number 1.2.
# Type it into a helper cell that's explicitly formatted as Text, then
# copy only the *value* (not the format) into D2 so D2 keeps the
# workbook's default (unstyled) cell format, matching the rest of the
# row.
$helper = $metaSheet.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "1.2"
$helper.Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()

$metaSheet.Range("E2").Value = "2021-09-02T07:51:49.703908Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:10.203102"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/254/?format=json"

# Re-select "data" as the active sheet (matches the original workbook,
# which opened with "data" as the only / first sheet).
$dataSheet.Activate()
